# SCRUM Update E, G, H und I
# Updates the "Backlog" (PO-Backlog) sheet with longer task descriptions and
# fills in the first Sprint-Backlog rows (incl. a new "Done" status column E).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Backlog" - expand the four short task descriptions
# ---------------------------------------------------------------------------
$backlog = $wb.Worksheets.Item(1)

$backlog.Range("C6").Value = "Umsetzung des Painters im Framework. Standard TicTactoe. Kreise und Kreuze und blaues Raster. Ausblick: karrierter Block als Hintergrund. Ausblick: Soundeffekte. Ausblick: Animationen."
$backlog.Range("C7").Value = "Umsetzung der Rules im Framework. Standard TicTacToe. Ausblick: Zeitlimit als Option."
$backlog.Range("C8").Value = "Umsetzung eines menschlichen Players im Framework. Standard TicTacToe mit Mausbedingung."
$backlog.Range("C9").Value = "Umsetzung eines PC gesteuerten Players im Framework. Standard TicTacToe per Zufall. Ausblick: Schwierigkeitsgrad."

# Move the view / selection to row 8 (whole row) like in the edited file.
$backlog.Activate()
$backlog.Rows.Item(8).Select()

# ---------------------------------------------------------------------------
# Sheet 2: "Sprint-Backlog" - add focus factor, due/done markers and fill in
# the first three sprint rows with data copied from the PO backlog.
# ---------------------------------------------------------------------------
$sprint = $wb.Worksheets.Item(2)
$sprint.Activate()

# Focus factor header cell.
$sprint.Range("D1").Value = "Focus Faktor: 0,5"

# "Done" marker column for the two already existing rows.
$sprint.Range("E3").Value = "Done"
$sprint.Range("E3").HorizontalAlignment = -4108
$sprint.Range("E3").VerticalAlignment = -4108

$sprint.Range("E4").Value = "Done"
$sprint.Range("E4").HorizontalAlignment = -4108
$sprint.Range("E4").VerticalAlignment = -4108

# Row 5: Spezifikation einer Umsetzungsidee ...
$sprint.Range("A5").Value = 2
$sprint.Range("B5").Value = "Spezifikation einer Umsetzungsidee für das Spiel TicTacToe"
$sprint.Range("C5").Value = "Spezifikation: Wie soll das TicTacToe aussehen? Wie soll das Aussehen technisch erreicht werden? Wie sollen Regeln umgesetzt werden? Soll es Sonderregeln geben, bzw. Sonderspielfelder, etc.? Welche Spieler (PC, ...) soll es geben? Etc.?"
$sprint.Range("D5").Value = "10min"
$sprint.Range("E5").Value = "Done"
$sprint.Range("E5").HorizontalAlignment = -4108
$sprint.Range("E5").VerticalAlignment = -4108

# Row 6: Implementierung der TicTacToe-Darstellung
$sprint.Range("A6").Value = 2
$sprint.Range("B6").Value = "Implementierung der TicTacToe-Darstellung"
$sprint.Range("C6").Value = "Umsetzung des Painters im Framework. Standard TicTactoe. Kreise und Kreuze und blaues Raster. Ausblick: karrierter Block als Hintergrund. Ausblick: Soundeffekte. Ausblick: Animationen."
$sprint.Range("D6").Value = "240min"
# touch E6 (stays empty) so the row keeps an (empty) cell in column E
$sprint.Range("E6").Font.Bold = $false

# Row 7: Implementierung eines menschlichen TicTacToe-Spielers
$sprint.Range("A7").Value = 2
$sprint.Range("B7").Value = "Implementierung eines menschlichen TicTacToe-Spielers"
$sprint.Range("C7").Value = "Umsetzung eines menschlichen Players im Framework. Standard TicTacToe mit Mausbedingung."
$sprint.Range("D7").Value = "120min"
# touch E7 (stays empty) so the row keeps an (empty) cell in column E
$sprint.Range("E7").Font.Bold = $false

# New selection / active cell on the Sprint-Backlog sheet.
$sprint.Range("E6").Select()
